$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: clear D4 (was "43063") back to an empty text cell, keep it text-typed
# (the quote-prefix trick forces Excel to treat the assignment as text, then we
# reset the style so no stray quotePrefix formatting sticks around).
$ws.Range("D4").Value = "'"
$ws.Range("D4").Style = "Normal"

# Row 4: numerator flips from TRUE to FALSE
$ws.Range("F4").Value = $false

# Row 4: numerator_desc text update
$ws.Range("G4").Value = "No screening recorded"

# Row 5: medicaid flips from FALSE to TRUE
$ws.Range("H5").Value = $true
